$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells keep their original text formatting
# (many values look numeric but must retain exact text, e.g. trailing zeros
# or "thousand.thousand.decimal" grouped notation), so force Text format
# before assigning values.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '65.926.92'
$ws.Range("E2").Value = '  -4.06%  '

$ws.Range("D3").Value = '3.295.59'
$ws.Range("E3").Value = '  -5.31%  '

$ws.Range("E4").Value = '  +0.22%  '

$ws.Range("D5").Value = '557.61'
$ws.Range("E5").Value = '  -4.10%  '

$ws.Range("D6").Value = '181.61'
$ws.Range("E6").Value = '  -4.21%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("E8").Value = '  -1.84%  '

$ws.Range("D9").Value = '3.290.87'
$ws.Range("E9").Value = '  -5.07%  '

$ws.Range("E10").Value = '  -7.07%  '

$ws.Range("E11").Value = '  -3.86%  '

$ws.Range("D12").Value = '47.75'
$ws.Range("E12").Value = '  -7.40%  '

$ws.Range("E13").Value = '  -6.06%  '

$ws.Range("D14").Value = '640.10'
$ws.Range("E14").Value = '  +0.78%  '

$ws.Range("D15").Value = '8.58'
$ws.Range("E15").Value = '  -5.48%  '

$ws.Range("D16").Value = '3.824.62'
$ws.Range("E16").Value = '  -4.71%  '

$ws.Range("D17").Value = '65.935.60'
$ws.Range("E17").Value = '  -4.20%  '

$ws.Range("D18").Value = '17.92'
$ws.Range("E18").Value = '  -0.82%  '

$ws.Range("E19").Value = '  -3.05%  '

$ws.Range("D20").Value = '3.294.60'
$ws.Range("E20").Value = '  -5.15%  '

$ws.Range("D21").Value = '11.44'
$ws.Range("E21").Value = '  -7.37%  '

$ws.Range("D22").Value = '0.907'
$ws.Range("E22").Value = '  -3.87%  '

$ws.Range("D23").Value = '17.86'
$ws.Range("E23").Value = '  +0.97%  '

$ws.Range("D24").Value = '107.69'
$ws.Range("E24").Value = '  +8.39%  '

$ws.Range("E25").Value = '  -6.33%  '

$ws.Range("E26").Value = '  -7.20%  '

$ws.Range("E27").Value = '  -6.07%  '

$ws.Range("D28").Value = '9.54'
$ws.Range("E28").Value = '  -4.80%  '

$ws.Range("D29").Value = '8.72'
$ws.Range("E29").Value = '  -4.51%  '

$ws.Range("D30").Value = '30.40'
$ws.Range("E30").Value = '  -6.16%  '

$ws.Range("D31").Value = '4.01'
$ws.Range("E31").Value = '  -1.02%  '

$ws.Range("D32").Value = '6.36'
$ws.Range("E32").Value = '  -4.95%  '

$ws.Range("D33").Value = '11.09'
$ws.Range("E33").Value = '  -4.05%  '

$ws.Range("D34").Value = '558.30'
$ws.Range("E34").Value = '  +11.91%  '

$ws.Range("E35").Value = '  -3.06%  '

$ws.Range("D36").Value = '57.30'
$ws.Range("E36").Value = '  -5.74%  '

$ws.Range("B37").Value = 'Dai'
$ws.Range("C37").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D37").Value = '0.999'

$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").Value = '3.685.15'
$ws.Range("E38").Value = '  -0.49%  '

$ws.Range("D39").Value = '3.50'
$ws.Range("E39").Value = '  -1.96%  '

$ws.Range("B40").Value = 'CoreDAO'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D40").Value = '3.59'
$ws.Range("E40").Value = '  +36.23%  '

$ws.Range("B41").Value = 'Fetch.AI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D41").Value = '2.74'
$ws.Range("E41").Value = '  -6.26%  '

$ws.Range("E42").Value = '  -9.52%  '

$ws.Range("E43").Value = '  -3.73%  '

$ws.Range("E44").Value = '  -6.32%  '

$ws.Range("D45").Value = '32.20'
$ws.Range("E45").Value = '  -5.83%  '

$ws.Range("D46").Value = '0.0417'
$ws.Range("E46").Value = '  -4.95%  '

$ws.Range("E47").Value = '  -2.45%  '

$ws.Range("E48").Value = '  -2.89%  '

$ws.Range("E49").Value = '  -5.84%  '

$ws.Range("D50").Value = '0.998'
$ws.Range("E50").Value = '  +0.10%  '

$ws.Range("D51").Value = '7.66'
